$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 33
$ws.Range("H33").Value = 124.5
$ws.Range("I33").Value = 91.85714
$ws.Range("J33").Value = 200.66667
$ws.Range("K33").Value = 91.85714
$ws.Range("L33").Value = 200.66667
$ws.Range("M33").Value = 137.14286
$ws.Range("N33").Value = -658.6666700000001

# Row 43
$ws.Range("H43").Value = 3750
$ws.Range("I43").Value = 3750
$ws.Range("J43").Value = 0
$ws.Range("K43").Value = 3750
$ws.Range("L43").Value = 0
$ws.Range("M43").Value = -3681
$ws.Range("N43").ClearContents()

# Row 51
$ws.Range("H51").Value = 2884.1428
$ws.Range("I51").Value = 2833
$ws.Range("J51").Value = 2922.5
$ws.Range("K51").Value = 2833
$ws.Range("L51").Value = 2922.5
$ws.Range("M51").Value = -2349
$ws.Range("N51").Value = -3890.5

# Row 80
$ws.Range("H80").Value = 1042
$ws.Range("I80").Value = 813
$ws.Range("J80").Value = 1500
$ws.Range("K80").Value = 2439
$ws.Range("L80").Value = 4500
$ws.Range("M80").Value = -1441
$ws.Range("N80").Value = -6496

# Row 83
$ws.Range("H83").Value = 1042
$ws.Range("I83").Value = 813
$ws.Range("J83").Value = 1500
$ws.Range("K83").Value = 7317
$ws.Range("L83").Value = 13500
$ws.Range("M83").Value = -2325
$ws.Range("N83").Value = -23484

# Row 129
$ws.Range("H129").Value = 1747
$ws.Range("J129").Value = 2493.6
$ws.Range("L129").Value = 7480.799999999999
$ws.Range("N129").Value = -17480.8

# Row 131
$ws.Range("H131").Value = 4052.1428
$ws.Range("I131").Value = 4060.8333
$ws.Range("K131").Value = 12182.4999
$ws.Range("M131").Value = -7142.499899999999

$ws = $wb.Worksheets.Item("ARM")
# Row 45
$ws.Range("H45").Value = 3265.5715
$ws.Range("I45").Value = 2328.3333
$ws.Range("J45").Value = 3968.5
$ws.Range("K45").Value = 2328.3333
$ws.Range("L45").Value = 3968.5
$ws.Range("M45").Value = -1951.3333
$ws.Range("N45").Value = -4722.5

# Row 61
$ws.Range("H61").Value = 1797.8
$ws.Range("I61").Value = 1797.8
$ws.Range("K61").Value = 1797.8
$ws.Range("M61").Value = -1585.8

# Row 74
$ws.Range("H74").Value = 5774.2085
$ws.Range("I74").Value = 5514.3
$ws.Range("K74").Value = 5514.3
$ws.Range("M74").Value = -4640.3

# Row 77
$ws.Range("H77").Value = 5774.2085
$ws.Range("I77").Value = 5514.3
$ws.Range("K77").Value = 27571.5
$ws.Range("M77").Value = -23203.5

# Row 102
$ws.Range("H102").Value = 6726.25
$ws.Range("I102").Value = 2936.6667
$ws.Range("K102").Value = 2936.6667
$ws.Range("M102").Value = -1314.6667

# Row 110
$ws.Range("H110").Value = 2335.0715
$ws.Range("I110").Value = 1720.75
$ws.Range("J110").Value = 3154.1667
$ws.Range("K110").Value = 1720.75
$ws.Range("L110").Value = 3154.1667
$ws.Range("M110").Value = 324.25
$ws.Range("N110").Value = -7244.1667

# Row 122
$ws.Range("H122").Value = 0
$ws.Range("I122").Value = 0
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 0
$ws.Range("L122").Value = 0
$ws.Range("M122").ClearContents()
$ws.Range("N122").ClearContents()

# Row 136
$ws.Range("H136").Value = 1797.8
$ws.Range("I136").Value = 1797.8
$ws.Range("K136").Value = 5393.4
$ws.Range("M136").Value = -2843.4

$ws = $wb.Worksheets.Item("BSM")
# Row 7
$ws.Range("H7").Value = 10341827
$ws.Range("I7").Value = 11500041
$ws.Range("K7").Value = 11500041
$ws.Range("M7").Value = -11499928

# Row 9
$ws.Range("H9").Value = 19950
$ws.Range("J9").Value = 19950
$ws.Range("L9").Value = 19950
$ws.Range("N9").Value = -20286

$ws = $wb.Worksheets.Item("CRP")
# Row 47
$ws.Range("H47").Value = 37575
$ws.Range("J47").Value = 37575
$ws.Range("L47").Value = 37575
$ws.Range("N47").Value = -38707

# Row 94
$ws.Range("H94").Value = 4768.875
$ws.Range("I94").Value = 1299
$ws.Range("K94").Value = 1299
$ws.Range("M94").Value = -848

# Row 122
$ws.Range("H122").Value = 787.86664
$ws.Range("I122").Value = 716.5
$ws.Range("J122").Value = 930.6
$ws.Range("K122").Value = 2149.5
$ws.Range("L122").Value = 2791.8
$ws.Range("M122").Value = 300.5
$ws.Range("N122").Value = -7691.8

$ws = $wb.Worksheets.Item("CUL")
# Row 11
$ws.Range("H11").Value = 333386.66
$ws.Range("I11").Value = 333386.66
$ws.Range("K11").Value = 1000159.98
$ws.Range("M11").Value = -1000019.98

# Row 34
$ws.Range("H34").Value = 2181.3845
$ws.Range("I34").Value = 671.375
$ws.Range("K34").Value = 2014.125
$ws.Range("M34").Value = -1930.125

# Row 86
$ws.Range("H86").Value = 343.33334
$ws.Range("I86").Value = 125
$ws.Range("J86").Value = 780
$ws.Range("K86").Value = 375
$ws.Range("L86").Value = 2340
$ws.Range("M86").Value = 811
$ws.Range("N86").Value = -4712

# Row 89
$ws.Range("H89").Value = 343.33334
$ws.Range("I89").Value = 125
$ws.Range("J89").Value = 780
$ws.Range("K89").Value = 1125
$ws.Range("L89").Value = 7020
$ws.Range("M89").Value = 4803
$ws.Range("N89").Value = -18876

# Row 134
$ws.Range("H134").Value = 3483.111
$ws.Range("I134").Value = 3483.111
$ws.Range("K134").Value = 10449.333
$ws.Range("M134").Value = -5379.332999999999

# Row 137
$ws.Range("H137").Value = 5133.3335
$ws.Range("I137").Value = 5000
$ws.Range("J137").Value = 5200
$ws.Range("K137").Value = 15000
$ws.Range("L137").Value = 15600
$ws.Range("M137").Value = -9900
$ws.Range("N137").Value = -25800

# Row 140
$ws.Range("H140").Value = 4095.182
$ws.Range("I140").Value = 3838.5557
$ws.Range("K140").Value = 11515.6671
$ws.Range("M140").Value = -6335.667099999999

$ws = $wb.Worksheets.Item("GSM")
# Row 4
$ws.Range("H4").Value = 380
$ws.Range("J4").Value = 380
$ws.Range("L4").Value = 380
$ws.Range("N4").Value = -604

# Row 23
$ws.Range("H23").Value = 748.3333
$ws.Range("J23").Value = 748.3333
$ws.Range("L23").Value = 748.3333
$ws.Range("N23").Value = -1194.3333

# Row 33
$ws.Range("H33").Value = 7874.75
$ws.Range("I33").Value = 7500
$ws.Range("J33").Value = 7999.6665
$ws.Range("K33").Value = 7500
$ws.Range("L33").Value = 7999.6665
$ws.Range("M33").Value = -7248
$ws.Range("N33").Value = -8503.666499999999

# Row 36
$ws.Range("H36").Value = 0
$ws.Range("J36").Value = 0
$ws.Range("L36").Value = 0
$ws.Range("N36").ClearContents()

# Row 46
$ws.Range("H46").Value = 10000
$ws.Range("I46").Value = 10000
$ws.Range("J46").Value = 0
$ws.Range("K46").Value = 10000
$ws.Range("L46").Value = 0
$ws.Range("M46").Value = -9844
$ws.Range("N46").ClearContents()

# Row 80
$ws.Range("H80").Value = 2374.25
$ws.Range("J80").Value = 2298
$ws.Range("L80").Value = 2298
$ws.Range("N80").Value = -4294

# Row 83
$ws.Range("H83").Value = 2374.25
$ws.Range("J83").Value = 2298
$ws.Range("L83").Value = 11490
$ws.Range("N83").Value = -21474

$ws = $wb.Worksheets.Item("LTW")
# Row 46
$ws.Range("H46").Value = 4578.5
$ws.Range("I46").Value = 3511.4
$ws.Range("J46").Value = 5063.5454
$ws.Range("K46").Value = 3511.4
$ws.Range("L46").Value = 5063.5454
$ws.Range("M46").Value = -3323.4
$ws.Range("N46").Value = -5439.5454

# Row 122
$ws.Range("H122").Value = 2474.5
$ws.Range("I122").Value = 2474.5
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 7423.5
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -4973.5
$ws.Range("N122").ClearContents()

# Row 132
$ws.Range("H132").Value = 5443.636
$ws.Range("I132").Value = 2328.125
$ws.Range("J132").Value = 13751.667
$ws.Range("K132").Value = 6984.375
$ws.Range("L132").Value = 41255.001
$ws.Range("M132").Value = -4454.375
$ws.Range("N132").Value = -46315.001

$ws = $wb.Worksheets.Item("WVR")
# Row 101
$ws.Range("H101").Value = 16380.4
$ws.Range("J101").Value = 16380.4
$ws.Range("L101").Value = 16380.4
$ws.Range("N101").Value = -22870.4

# Row 133
$ws.Range("H133").Value = 0
$ws.Range("J133").Value = 0
$ws.Range("L133").Value = 0
$ws.Range("N133").ClearContents()
